# feat: add 2022-Q3 data
#
# - Inserts a new "2022-Q3" worksheet (fund holdings for that quarter),
#   positioned right before the existing "2022-Q1" tab.
# - Updates the "总计" (totals) summary sheet with a new leading row for
#   2022-Q3, shifting the previously-first rows (2022-Q1 / 2021-Q1) down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Build the "2022-Q3" sheet by duplicating "2022-Q1" (this keeps all
#    of its header/column formatting intact) and then rewriting its
#    cell contents in place.
# ---------------------------------------------------------------------
$q1Sheet = $wb.Worksheets.Item("2022-Q1")
$q1Sheet.Copy($q1Sheet)                      # inserts the copy right before 2022-Q1
$q3Sheet = $wb.Worksheets.Item("2022-Q1 (2)")
$q3Sheet.Name = "2022-Q3"

# The source sheet has 4 data rows; 2022-Q3 only needs 2, so drop rows 4-5.
$q3Sheet.Rows.Item(4).Resize(2).Delete()

# Header row (unchanged text, but rewrite explicitly for clarity/safety)
$q3Sheet.Range("B1").Value = "基金代码"
$q3Sheet.Range("C1").Value = "基金名称"
$q3Sheet.Range("D1").Value = "基金规模"
$q3Sheet.Range("E1").Value = "股票总仓位"
$q3Sheet.Range("F1").Value = "仓位占比"
$q3Sheet.Range("G1").Value = "持有市值(亿元)"
$q3Sheet.Range("H1").Value = "仓位排名"

# Row 2 - 013132 / 创金合信文娱媒体股票型发起式证券投资基金A
$q3Sheet.Range("A2").Value = 0
$q3Sheet.Range("B2").Value = "'013132"
$q3Sheet.Range("C2").Value = "创金合信文娱媒体股票型发起式证券投资基金A"
$q3Sheet.Range("D2").Value = "'0.06"
$q3Sheet.Range("E2").Value = "'91.90"
$q3Sheet.Range("F2").Value = "'5.18"
$q3Sheet.Range("G2").Value = "'0.0031"
$q3Sheet.Range("H2").Value = 3

# Row 3 - 013133 / 创金合信文娱媒体股票型发起式证券投资基金C
$q3Sheet.Range("A3").Value = 1
$q3Sheet.Range("B3").Value = "'013133"
$q3Sheet.Range("C3").Value = "创金合信文娱媒体股票型发起式证券投资基金C"
$q3Sheet.Range("D3").Value = "'0.05"
$q3Sheet.Range("E3").Value = "'91.90"
$q3Sheet.Range("F3").Value = "'5.18"
$q3Sheet.Range("G3").Value = "'0.0026"
$q3Sheet.Range("H3").Value = 3

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: new first data row for 2022-Q3,
#    with the previously existing rows (2022-Q1 / 2021-Q1) re-written
#    below it and their running index bumped by one.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Resize(10).ClearContents()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.01

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q1"
$total.Range("C3").Value = 4
$total.Range("D3").Value = 0.13

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q1"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.01

# A4 is a brand-new cell; copy the bold/bordered/centered style already
# used by A2/A3 onto it so all three index cells stay visually consistent.
$total.Range("A3").Copy()
$total.Range("A4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

Write-Host "2022-Q3 sheet added; 总计 sheet updated"
